# Add two new "add tag" log rows to Sheet1, right after the existing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$timestamp = 1749897806

$ws.Range("A45").Value = $timestamp
$ws.Range("B45").Value = "add"
$ws.Range("C45").Value = "tag"
$ws.Range("D45").Value = "images_aeriennes_1"

$ws.Range("A46").Value = $timestamp
$ws.Range("B46").Value = "add"
$ws.Range("C46").Value = "tag"
$ws.Range("D46").Value = "images_aeriennes_2"
